$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Nifty": roll the date window forward by one trading day.
# Old G (Feb 15) drops out, old H (Feb 18) becomes new G, and a new H
# (Feb 19) is entered. All pivot / camarilla formulas downstream
# recalculate automatically.
# ---------------------------------------------------------------------
$wsNifty = $wb.Worksheets.Item("Nifty")

$wsNifty.Range("G1").Value = 43514
$wsNifty.Range("H1").Value = 43515

$wsNifty.Range("G2").Value = 10759.9
$wsNifty.Range("H2").Value = 10722.85

$wsNifty.Range("G3").Value = 10628.4
$wsNifty.Range("H3").Value = 10585.65

$wsNifty.Range("G4").Value = 10640.95
$wsNifty.Range("H4").Value = 10604.35

# Elliott wave levels (rows 38-48) - hard coded numbers, also rolled
# forward one column with several brand-new annotated levels.
$wsNifty.Range("G38").ClearContents()

$wsNifty.Range("G39").ClearContents()
$wsNifty.Range("H39").Value = 10810.5214

$wsNifty.Range("G40").Value = 10854.110999999999
$wsNifty.Range("H40").Value = 10759

$wsNifty.Range("G41").Value = 10810.5214
$wsNifty.Range("H41").Value = 10638.0604
$wsNifty.Range("H41").Font.Bold = $true

$wsNifty.Range("H42").Value = 10618.029199999999

$wsNifty.Range("H44").Value = 10583.35

$wsNifty.Range("H45").Value = 10569.7786
$wsNifty.Range("I45").NumberFormat = "0%"

$wsNifty.Range("H46").Value = 10550.428
$wsNifty.Range("I46").NumberFormat = "0%"

$wsNifty.Range("H47").Value = 10511.05

$wsNifty.Range("H48").Value = 10452.321399999999

# New annotation labels in column I describing each Elliott level
$wsNifty.Range("I44").Value = "C 100% Proj"
$wsNifty.Range("I46").Value = "C 123% Proj"
$wsNifty.Range("I45").Value = "5 38% Proj"
$wsNifty.Range("I47").Value = "5 50% Proj"
$wsNifty.Range("I48").Value = "5 60% Proj"
$wsNifty.Range("I42").Value = "C 23% Ret"
$wsNifty.Range("I41").Value = "C 38% Ret Bullish"
$wsNifty.Range("I40").Value = "Wave A End"
$wsNifty.Range("I39").Value = "38% Ret 0 to 3"

# ---------------------------------------------------------------------
# Sheet "Elliot": Start/End points for the Fibonacci retracement grids
# (columns H/J/L/N/P) roll forward the same way the Nifty price window
# did. All of the 23.6% / 38.2% / 50% / ... retracement + projection
# rows are formulas off these cells and recompute on their own.
# ---------------------------------------------------------------------
$wsElliot = $wb.Worksheets.Item("Elliot")

$wsElliot.Range("J6").ClearContents()
$wsElliot.Range("L6").Value = 10759.9
$wsElliot.Range("N6").Value = 10722.85
$wsElliot.Range("P6").ClearContents()

$wsElliot.Range("J9").ClearContents()
$wsElliot.Range("L9").Value = 10620.4
$wsElliot.Range("N9").Value = 10585.65
$wsElliot.Range("P9").ClearContents()

$wsElliot.Range("H12").Value = 10759.9
$wsElliot.Range("J12").ClearContents()
$wsElliot.Range("L12").Value = 10722.85
$wsElliot.Range("N12").ClearContents()

# ---------------------------------------------------------------------
# View / selection bookkeeping, applied last so "Nifty" ends up as the
# active (tabSelected) sheet again, matching the original workbook.
# ---------------------------------------------------------------------
$wsElliot.Activate()
$wsElliot.Range("H25:H27").Select()

$wsFib = $wb.Worksheets.Item("Fibonnacci")
$wsFib.Activate()
$wsFib.Range("B17").Select()

$wsNifty.Activate()
$wsNifty.Range("M39").Select()
